$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "user" columns (B:E) were empty -> fill in the new
#     "/refreshToken" route; "chat" columns (J:M) get their route path and
#     returned-value note updated ("/new" -> "/create", note expanded).
$ws.Range("B10").Value = "POST"
$ws.Range("C10").Value = "/refreshToken"
$ws.Range("D10").Value = "create new accessToken from refreshToken"
$ws.Range("E10").Value = "-"

$ws.Range("K10").Value = "/create"
$ws.Range("M10").Value = "returned value - new chat id, need authentication!"

# --- Row 11: "user" columns (B:E) were empty -> fill in the new "/data"
#     route (get user data according to the access token).
$ws.Range("B11").Value = "GET"
$ws.Range("C11").Value = "/data"
$ws.Range("D11").Value = "get user data according accessToken"
$ws.Range("E11").Value = "returned value - user object with username, user id and image URL"

# Both rows grow to fit the extra wrapped text (was default height).
$ws.Rows.Item(10).RowHeight = 42
$ws.Rows.Item(11).RowHeight = 42

# Selection/scroll position left where the author was last working.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E10").Select()
